$wb = $excel.ActiveWorkbook

$wsWebsite = $wb.Worksheets.Item("websiteRegistarion")
$wsWebsite.Range("A2").Value = "testweb30july22@gmail.com"

$wsKyc = $wb.Worksheets.Item("allreadyKyc")
$wsKyc.Range("A2").Value = "TEEPT2091J"
